$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "General Meeting Pizza and Soda" line item from Club Administration.
$ws.Rows.Item(35).Delete()

# Remove the "Food" line item from the SEDS Conference sub-section
# (originally row 43, now row 42 after the first deletion shifted rows up).
$ws.Rows.Item(42).Delete()

# Update Engineering budget figures.
$ws.Range("B17").Value = 100     # Launch Pad Components
$ws.Range("B18").Value = 1000    # Medium Power Rocket Components
$ws.Range("B19").Value = 3000    # Hybrid Rocket Hardware
$ws.Range("B20").Value = 500     # Miscellaneous Expenses

# Update the attendance note and travel costs for the SEDS Conference.
$ws.Range("A38").Value = "(16 Members in Attendence)"
$ws.Range("B39").Value = 5000    # Airplane Ticket and General Travel
$ws.Range("B40").Value = 1000    # Hotel Housing (3 Nights)

$ws.Range("I31").Select()
